# Update row 37 from "World Cup Qualifiers (UEFA)" to "SWPL"
# and remove rows 38-48 (old women's football leagues data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 37 with new competition data
# (set C37 first so the shared-string table order matches the target workbook)
$ws.Range("C37").Value = "ceu82myq9rpq841ts3jl7uvis"
$ws.Range("B37").Value = "SWPL"

# Delete rows 38 through 48 (shift cells up)
$ws.Range("A38:D48").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Update the view: scroll to top-left A4, selection at B38
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B38").Select()
